$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Wed Nov 02 16:24:38 EDT 2022"
$ws.Range("B3").Value = "Wed Nov 02 16:24:49 EDT 2022"
$ws.Range("B4").Value = "Wed Nov 02 16:25:01 EDT 2022"
